# Insert a new data row at row 84 on the single worksheet.
# This pushes the existing rows 84-131 down to 85-132 (preserving their
# values/formatting untouched) and fills the freshly inserted row 84 with
# a new "Arveja Verde" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 84..131 down by one row.
$ws.Rows(84).Insert()

# Populate the newly inserted row 84.
$ws.Range("A84").Value = 9
$ws.Range("B84").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C84").Value = "Metropolitana"
$ws.Range("D84").Value = 44806
$ws.Range("E84").Value = 13
$ws.Range("F84").Value = 100112022
$ws.Range("G84").Value = "Arveja Verde"
$ws.Range("H84").Value = "Perfection"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 25
$ws.Range("K84").Value = 38000
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = 38800
$ws.Range("N84").Value = "`$/malla 25 kilos"
$ws.Range("O84").Value = "Provincia de Huasco"
$ws.Range("P84").Value = 1552
$ws.Range("Q84").Value = 25
$ws.Range("R84").Value = "Hortaliza"
